$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: proper column headers ---
# (previously this row just duplicated row 2's data instead of holding labels)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Give the newly added header cells (H1:N1) the same look as the rest of the
# header row (bold, centered, thin border) by copying B1's formatting.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: append the metadata columns every other sheet already has ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("K2").Value = "何欣純"
$ws.Range("L2").Value = 1733
$ws.Range("M2").Value = "tmp2e891"
$ws.Range("N2").Value = 32

# "date" (J2) must stay plain text ("2012-04-30"), not get auto-converted to a
# date serial. Build it as text in an unused helper cell, then copy just the
# value across so no stray number-format sneaks onto J2.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "2012-04-30"
$helper.Copy()
$ws.Range("J2").PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = $false

# Match the plain data-row look (no border/bold) for the new row-2 cells.
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
